$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price column (D) cells to be treated as text so values like
# "1.0000" / "0.9993" / "29.140.04" are not auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '29.140.04'
$ws.Range("E2").Value = '  -2.00%  '

# Row 3
$ws.Range("D3").Value = '1.839.81'
$ws.Range("E3").Value = '  -1.35%  '

# Row 4
$ws.Range("D4").Value = '0.9993'
$ws.Range("E4").Value = '  -0.09%  '

# Row 5
$ws.Range("D5").Value = '239.90'
$ws.Range("E5").Value = '  -2.69%  '

# Row 6
$ws.Range("D6").Value = '0.6813'

# Row 7
$ws.Range("D7").Value = '1.0000'
$ws.Range("E7").Value = '  -0.07%  '

# Row 8
$ws.Range("D8").Value = '0.2993'
$ws.Range("E8").Value = '  -2.96%  '

# Row 9
$ws.Range("D9").Value = '0.07459'
$ws.Range("E9").Value = '  -4.08%  '

# Row 10
$ws.Range("D10").Value = '23.27'
$ws.Range("E10").Value = '  -2.23%  '

# Row 11
$ws.Range("D11").Value = '0.07645'
$ws.Range("E11").Value = '  -2.53%  '

# Row 12
$ws.Range("D12").Value = '1.840.18'
$ws.Range("E12").Value = '  -1.34%  '

# Row 13
$ws.Range("D13").Value = '5.034'
$ws.Range("E13").Value = '  -2.88%  '

# Row 14
$ws.Range("D14").Value = '0.6808'
$ws.Range("E14").Value = '  -2.23%  '

# Row 15
$ws.Range("E15").Value = '  -5.82%  '

# Row 16
$ws.Range("D16").Value = '6.153'
$ws.Range("E16").Value = '  -7.41%  '

# Row 17
$ws.Range("D17").Value = '29.132.95'
$ws.Range("E17").Value = '  -2.05%  '

# Row 18
$ws.Range("D18").Value = '0.000008215'
$ws.Range("E18").Value = '  -2.16%  '

# Row 19
$ws.Range("D19").Value = '2.079.85'
$ws.Range("E19").Value = '  -2.02%  '

# Row 20
$ws.Range("D20").Value = '230.65'
$ws.Range("E20").Value = '  -5.58%  '

# Row 21
$ws.Range("D21").Value = '12.51'
$ws.Range("E21").Value = '  -2.53%  '

# Row 22
$ws.Range("D22").Value = '1.0000'
$ws.Range("E22").Value = '  -0.05%  '

# Row 23
$ws.Range("D23").Value = '7.335'
$ws.Range("E23").Value = '  -4.39%  '

# Row 24
$ws.Range("D24").Value = '0.9998'
$ws.Range("E24").Value = '  -0.08%  '

# Row 25
$ws.Range("D25").Value = '161.17'
$ws.Range("E25").Value = '  +0.43%  '

# Row 26
$ws.Range("D26").Value = '0.1429'
$ws.Range("E26").Value = '  -6.03%  '

# Row 27
$ws.Range("D27").Value = '8.704'
$ws.Range("E27").Value = '  -3.07%  '

# Row 28
$ws.Range("D28").Value = '18.07'
$ws.Range("E28").Value = '  -1.98%  '

# Row 29
$ws.Range("E29").Value = '  -3.21%  '

# Row 30
$ws.Range("D30").Value = '4.260'
$ws.Range("E30").Value = '  -0.55%  '

# Row 31
$ws.Range("D31").Value = '4.139'
$ws.Range("E31").Value = '  -1.71%  '

# Row 32
$ws.Range("D32").Value = '1.196'
$ws.Range("E32").Value = '  -0.31%  '

# Row 33
$ws.Range("D33").Value = '0.05376'
$ws.Range("E33").Value = '  +5.15%  '

# Row 34
$ws.Range("D34").Value = '0.7544'
$ws.Range("E34").Value = '  -4.62%  '

# Row 35
$ws.Range("D35").Value = '1.849'
$ws.Range("E35").Value = '  -3.89%  '

# Row 36
$ws.Range("E36").Value = '  -2.43%  '

# Row 37
$ws.Range("D37").Value = '2.687'
$ws.Range("E37").Value = '  -0.32%  '

# Row 38
$ws.Range("D38").Value = '1.311.16'
$ws.Range("E38").Value = '  -2.56%  '

# Row 39
$ws.Range("D39").Value = '0.01830'
$ws.Range("E39").Value = '  -3.17%  '

# Row 40
$ws.Range("D40").Value = '2.723'
$ws.Range("E40").Value = '  -0.66%  '

# Row 41
$ws.Range("D41").Value = '0.9430'
$ws.Range("E41").Value = '  -2.25%  '

# Row 42
$ws.Range("D42").Value = '6.046'
$ws.Range("E42").Value = '  -0.05%  '

# Row 43
$ws.Range("D43").Value = '104.90'
$ws.Range("E43").Value = '  -1.76%  '

# Row 45
$ws.Range("B45").Value = 'XinFinNetwork'
$ws.Range("C45").Value = 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'
$ws.Range("D45").Value = '0.07850'
$ws.Range("E45").Value = '  +22.70%  '

# Row 46
$ws.Range("B46").Value = 'RocketPoolETH'
$ws.Range("C46").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D46").Value = '1.981.05'
$ws.Range("E46").Value = '  -1.64%  '

# Row 47
$ws.Range("D47").Value = '0.5179'
$ws.Range("E47").Value = '  -0.34%  '

# Row 48
$ws.Range("E48").Value = '  -2.16%  '

# Row 49
$ws.Range("D49").Value = '64.22'
$ws.Range("E49").Value = '  -2.00%  '

# Row 50
$ws.Range("D50").Value = '1.775'
$ws.Range("E50").Value = '  -1.13%  '

# Row 51
$ws.Range("D51").Value = '9.449'
$ws.Range("E51").Value = '  -3.67%  '
